$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 for the Driver Tracking module task
$ws.Range("A12").Value = "#100010"

# Fix typo in existing cell B11: "CreateMain Page" -> "Create Main Page"
$ws.Range("B11").Value = "Create Main Page"

$ws.Range("B12").Value = "Create Driver Tracking module"

# Move the active selection cursor as recorded in the edit
$ws.Range("B20").Select()
